$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("A1").Value = "SNO"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Pass"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Test test 2"
$ws.Range("C2").Value = "Pass"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Test test 2"
$ws.Range("C3").Value = "Pass"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Test test 2"
$ws.Range("C4").Value = "Pass"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Test test 2"
$ws.Range("C5").Value = "Pass"

$ws.Range("A6").Value = 21

# Update selection to match the target state (active cell A5 within A4:A6)
$ws.Range("A4:A6").Select()
$excel.ActiveWindow.ActiveCell = $ws.Range("A5")
